# Scheduled-runner style market/profit refresh across the Leve-profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for a
# handful of rows on each sheet, reflecting newer market-board snapshots.

$wb = $excel.ActiveWorkbook

# ---- ALC ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value  = 645.75
$ws.Range("I15").Value  = 645.75
$ws.Range("K15").Value  = 1937.25
$ws.Range("M15").Value  = -1768.25

$ws.Range("H133").Value = 34500
$ws.Range("J133").Value = 34500
$ws.Range("L133").Value = 34500
$ws.Range("N133").Value = -44620

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 34992
$ws.Range("J136").Value = 34992
$ws.Range("L136").Value = 34992
$ws.Range("N136").Value = -45192

$ws.Range("H140").Value = 76036.664
$ws.Range("J140").Value = 73847.5
$ws.Range("L140").Value = 73847.5
$ws.Range("N140").Value = -84207.5

# ---- ARM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value  = 12465.797
$ws.Range("I32").Value  = 13068.266
$ws.Range("J32").Value  = 8610
$ws.Range("K32").Value  = 13068.266
$ws.Range("L32").Value  = 8610
$ws.Range("M32").Value  = -12781.266
$ws.Range("N32").Value  = -9184

$ws.Range("H123").Value = 55428.5
$ws.Range("J123").Value = 55428.5
$ws.Range("L123").Value = 55428.5
$ws.Range("N123").Value = -65228.5

$ws.Range("H134").Value = 66220
$ws.Range("J134").Value = 66220
$ws.Range("L134").Value = 66220
$ws.Range("N134").Value = -76360

$ws.Range("H138").Value = 56979.715
$ws.Range("J138").Value = 56979.715
$ws.Range("L138").Value = 56979.715
$ws.Range("N138").Value = -67259.715

$ws.Range("H139").Value = 89031.89
$ws.Range("J139").Value = 89031.89
$ws.Range("L139").Value = 89031.89
$ws.Range("N139").Value = -99311.89

# ---- BSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5071.5
$ws.Range("I105").Value = 5071.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5071.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3324.5
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 75647.914
$ws.Range("J132").Value = 75647.914
$ws.Range("L132").Value = 75647.914
$ws.Range("N132").Value = -85767.914

$ws.Range("H135").Value = 56282.312
$ws.Range("J135").Value = 56282.312
$ws.Range("L135").Value = 56282.312
$ws.Range("N135").Value = -66422.31200000001

$ws.Range("H140").Value = 37640.777
$ws.Range("J140").Value = 37640.777
$ws.Range("L140").Value = 37640.777
$ws.Range("N140").Value = -48000.777

# ---- CRP -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 76579
$ws.Range("J135").Value = 133158
$ws.Range("L135").Value = 133158
$ws.Range("N135").Value = -143298

$ws.Range("H137").Value = 54542.223
$ws.Range("J137").Value = 74176
$ws.Range("L137").Value = 74176
$ws.Range("N137").Value = -84376

$ws.Range("H138").Value = 52397
$ws.Range("J138").Value = 52397
$ws.Range("L138").Value = 52397
$ws.Range("N138").Value = -62677

$ws.Range("H140").Value = 71100
$ws.Range("J140").Value = 71100
$ws.Range("L140").Value = 71100
$ws.Range("N140").Value = -81460

# ---- CUL -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2389
$ws.Range("I103").Value = 625
$ws.Range("J103").Value = 3271
$ws.Range("K103").Value = 1875
$ws.Range("L103").Value = 9813
$ws.Range("M103").Value = -996
$ws.Range("N103").Value = -11571

# ---- GSM -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2201.9
$ws.Range("I132").Value = 1770.9615
$ws.Range("J132").Value = 5003
$ws.Range("K132").Value = 5312.8845
$ws.Range("L132").Value = 15009
$ws.Range("M132").Value = -2782.8845
$ws.Range("N132").Value = -20069

$ws.Range("H133").Value = 59826.668
$ws.Range("J133").Value = 59826.668
$ws.Range("L133").Value = 59826.668
$ws.Range("N133").Value = -69946.66800000001

$ws.Range("H135").Value = 48737.2
$ws.Range("J135").Value = 48737.2
$ws.Range("L135").Value = 48737.2
$ws.Range("N135").Value = -58877.2

$ws.Range("H138").Value = 54233.332
$ws.Range("J138").Value = 54233.332
$ws.Range("L138").Value = 54233.332
$ws.Range("N138").Value = -64513.332

$ws.Range("H139").Value = 167469
$ws.Range("J139").Value = 167469
$ws.Range("L139").Value = 167469
$ws.Range("N139").Value = -177749

$ws.Range("H140").Value = 48987.668
$ws.Range("J140").Value = 48987.668
$ws.Range("L140").Value = 48987.668
$ws.Range("N140").Value = -59347.668

$ws.Range("H141").Value = 41196
$ws.Range("J141").Value = 41196
$ws.Range("L141").Value = 41196
$ws.Range("N141").Value = -51556

# ---- LTW -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 31313
$ws.Range("J108").Value = 31313
$ws.Range("L108").Value = 31313
$ws.Range("N108").Value = -38993

$ws.Range("H132").Value = 4513.6816
$ws.Range("I132").Value = 4719.6772
$ws.Range("K132").Value = 14159.0316
$ws.Range("M132").Value = -11629.0316

$ws.Range("H134").Value = 74957.25
$ws.Range("J134").Value = 74957.25
$ws.Range("L134").Value = 74957.25
$ws.Range("N134").Value = -85097.25

$ws.Range("H135").Value = 115189.91
$ws.Range("J135").Value = 115189.91
$ws.Range("L135").Value = 115189.91
$ws.Range("N135").Value = -125329.91

$ws.Range("H137").Value = 73097.2
$ws.Range("I137").Value = 20390
$ws.Range("J137").Value = 78953.55499999999
$ws.Range("K137").Value = 20390
$ws.Range("L137").Value = 78953.55499999999
$ws.Range("M137").Value = -15290
$ws.Range("N137").Value = -89153.55499999999

$ws.Range("H139").Value = 37205.668
$ws.Range("J139").Value = 37205.668
$ws.Range("L139").Value = 37205.668
$ws.Range("N139").Value = -47485.668

$ws.Range("H140").Value = 69285.75
$ws.Range("J140").Value = 69285.75
$ws.Range("L140").Value = 69285.75
$ws.Range("N140").Value = -79645.75

$ws.Range("H141").Value = 45285.832
$ws.Range("J141").Value = 45285.832
$ws.Range("L141").Value = 45285.832
$ws.Range("N141").Value = -55645.832

# ---- WVR -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2155.8462
$ws.Range("I132").Value = 1626.4117
$ws.Range("J132").Value = 3155.889
$ws.Range("K132").Value = 4879.2351
$ws.Range("L132").Value = 9467.667000000001
$ws.Range("M132").Value = -2349.2351
$ws.Range("N132").Value = -14527.667

$ws.Range("H133").Value = 42372.2
$ws.Range("J133").Value = 42372.2
$ws.Range("L133").Value = 42372.2
$ws.Range("N133").Value = -52492.2

$ws.Range("H135").Value = 58160.453
$ws.Range("J135").Value = 58160.453
$ws.Range("L135").Value = 58160.453
$ws.Range("N135").Value = -68300.45300000001

$ws.Range("H136").Value = 1289.9
$ws.Range("I136").Value = 1204.2363
$ws.Range("J136").Value = 1604
$ws.Range("K136").Value = 3612.7089
$ws.Range("L136").Value = 4812
$ws.Range("M136").Value = -1062.7089
$ws.Range("N136").Value = -9912

$ws.Range("H139").Value = 57707.145
$ws.Range("J139").Value = 57707.145
$ws.Range("L139").Value = 57707.145
$ws.Range("N139").Value = -67987.14499999999

$ws.Range("H140").Value = 34993.332
$ws.Range("J140").Value = 34993.332
$ws.Range("L140").Value = 34993.332
$ws.Range("N140").Value = -45353.332

$ws.Range("H141").Value = 83159.91
$ws.Range("J141").Value = 83159.91
$ws.Range("L141").Value = 83159.91
$ws.Range("N141").Value = -93519.91
